$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.289.73'
$ws.Range('E2').Value = '  -3.00%  '

# Row 3
$ws.Range('D3').Value = '1.854.54'
$ws.Range('E3').Value = '  -3.74%  '

# Row 4
$ws.Range('E4').Value = '  -0.16%  '

# Row 5
$ws.Range('D5').Value = "'323.73"
$ws.Range('E5').Value = '  -1.83%  '

# Row 6
$ws.Range('E6').Value = '  -0.15%  '

# Row 7
$ws.Range('D7').Value = "'0.4536"
$ws.Range('E7').Value = '  -4.06%  '

# Row 8
$ws.Range('D8').Value = "'0.3864"
$ws.Range('E8').Value = '  -4.97%  '

# Row 9
$ws.Range('D9').Value = "'48.36"
$ws.Range('E9').Value = '  -8.82%  '

# Row 10
$ws.Range('D10').Value = "'0.07918"
$ws.Range('E10').Value = '  -6.16%  '

# Row 11
$ws.Range('D11').Value = "'1.014"
$ws.Range('E11').Value = '  -3.17%  '

# Row 12
$ws.Range('D12').Value = "'21.36"
$ws.Range('E12').Value = '  -4.04%  '

# Row 13
$ws.Range('D13').Value = '1.862.30'
$ws.Range('E13').Value = '  -4.06%  '

# Row 14
$ws.Range('D14').Value = "'5.901"
$ws.Range('E14').Value = '  -3.26%  '

# Row 15
$ws.Range('D15').Value = "'7.139"
$ws.Range('E15').Value = '  -5.06%  '

# Row 16
$ws.Range('E16').Value = '  -0.31%  '

# Row 17
$ws.Range('D17').Value = "'85.84"
$ws.Range('E17').Value = '  -5.32%  '

# Row 18
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = "'0.06551"
$ws.Range('E18').Value = '  -0.47%  '

# Row 19
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = "'0.00001022"
$ws.Range('E19').Value = '  -4.10%  '

# Row 20
$ws.Range('D20').Value = "'17.04"
$ws.Range('E20').Value = '  -6.23%  '

# Row 21
$ws.Range('E21').Value = '  -0.15%  '

# Row 22
$ws.Range('D22').Value = "'5.500"
$ws.Range('E22').Value = '  -4.59%  '

# Row 23
$ws.Range('D23').Value = '27.296.02'
$ws.Range('E23').Value = '  -3.08%  '

# Row 24
$ws.Range('D24').Value = "'10.87"
$ws.Range('E24').Value = '  -4.91%  '

# Row 25
$ws.Range('E25').Value = '  +0.28%  '

# Row 26
$ws.Range('D26').Value = '2.078.55'
$ws.Range('E26').Value = '  -4.30%  '

# Row 27
$ws.Range('D27').Value = "'153.65"
$ws.Range('E27').Value = '  -0.43%  '

# Row 28
$ws.Range('D28').Value = "'19.86"
$ws.Range('E28').Value = '  -1.35%  '

# Row 29
$ws.Range('D29').Value = "'2.058"
$ws.Range('E29').Value = '  -4.49%  '

# Row 30
$ws.Range('D30').Value = "'5.466"
$ws.Range('E30').Value = '  -4.78%  '

# Row 31
$ws.Range('D31').Value = "'121.21"
$ws.Range('E31').Value = '  -2.06%  '

# Row 32
$ws.Range('D32').Value = "'0.09314"

# Row 33
$ws.Range('D33').Value = "'0.9344"
$ws.Range('E33').Value = '  -4.04%  '

# Row 34
$ws.Range('D34').Value = "'1.458"
$ws.Range('E34').Value = '  +1.19%  '

# Row 35
$ws.Range('D35').Value = "'3.585"
$ws.Range('E35').Value = '  -1.45%  '

# Row 36
$ws.Range('D36').Value = "'5.265"
$ws.Range('E36').Value = '  -5.29%  '

# Row 37
$ws.Range('D37').Value = "'0.02223"
$ws.Range('E37').Value = '  -3.97%  '

# Row 38
$ws.Range('D38').Value = "'0.05996"
$ws.Range('E38').Value = '  -2.87%  '

# Row 39
$ws.Range('D39').Value = "'1.220"
$ws.Range('E39').Value = '  -1.45%  '

# Row 40
$ws.Range('D40').Value = "'8.053"
$ws.Range('E40').Value = '  -10.86%  '

# Row 41
$ws.Range('D41').Value = "'1.001"
$ws.Range('E41').Value = '  -0.14%  '

# Row 42
$ws.Range('D42').Value = "'0.5907"
$ws.Range('E42').Value = '  -4.37%  '

# Row 43
$ws.Range('D43').Value = "'0.1881"
$ws.Range('E43').Value = '  -1.21%  '

# Row 44
$ws.Range('D44').Value = "'10.13"
$ws.Range('E44').Value = '  -8.42%  '

# Row 45
$ws.Range('E45').Value = '  -1.16%  '

# Row 46
$ws.Range('D46').Value = "'0.5618"
$ws.Range('E46').Value = '  -4.79%  '

# Row 47
$ws.Range('D47').Value = "'12.05"
$ws.Range('E47').Value = '  -5.75%  '

# Row 48
$ws.Range('E48').Value = '  -2.85%  '

# Row 49
$ws.Range('D49').Value = "'1.915"
$ws.Range('E49').Value = '  -6.04%  '

# Row 50
$ws.Range('D50').Value = "'0.06736"
$ws.Range('E50').Value = '  -1.26%  '

# Row 51
$ws.Range('D51').Value = "'108.63"
$ws.Range('E51').Value = '  -1.29%  '
